$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3168.0667
$ws.Range("J17").Value = 3318.1667
$ws.Range("L17").Value = 9954.500100000001
$ws.Range("N17").Value = -10290.5001
$ws.Range("H40").Value = 4301.5713
$ws.Range("J40").Value = 4633.88
$ws.Range("L40").Value = 4633.88
$ws.Range("N40").Value = -4983.88
$ws.Range("H55").Value = 176.64706
$ws.Range("I55").Value = 126.57143
$ws.Range("J55").Value = 211.7
$ws.Range("K55").Value = 126.57143
$ws.Range("L55").Value = 211.7
$ws.Range("M55").Value = 87.42856999999999
$ws.Range("N55").Value = -639.7
$ws.Range("H80").Value = 2456.0264
$ws.Range("I80").Value = 634.1053000000001
$ws.Range("J80").Value = 4277.9473
$ws.Range("K80").Value = 1902.3159
$ws.Range("L80").Value = 12833.8419
$ws.Range("M80").Value = -904.3159000000001
$ws.Range("N80").Value = -14829.8419
$ws.Range("H83").Value = 2456.0264
$ws.Range("I83").Value = 634.1053000000001
$ws.Range("J83").Value = 4277.9473
$ws.Range("K83").Value = 5706.947700000001
$ws.Range("L83").Value = 38501.5257
$ws.Range("M83").Value = -714.9477000000006
$ws.Range("N83").Value = -48485.5257
$ws.Range("H112").Value = 7044.222
$ws.Range("J112").Value = 7044.222
$ws.Range("L112").Value = 21132.666
$ws.Range("N112").Value = -23348.666
$ws.Range("H125").Value = 7578267
$ws.Range("I125").Value = 1387.75
$ws.Range("J125").Value = 10419597
$ws.Range("K125").Value = 12489.75
$ws.Range("L125").Value = 93776373
$ws.Range("M125").Value = -10029.75
$ws.Range("N125").Value = -93781293
$ws.Range("H132").Value = 20410682
$ws.Range("I132").Value = 22224818
$ws.Range("J132").Value = 1652.25
$ws.Range("K132").Value = 66674454
$ws.Range("L132").Value = 4956.75
$ws.Range("M132").Value = -66671924
$ws.Range("N132").Value = -10016.75
$ws.Range("H137").Value = 66478.39
$ws.Range("I137").Value = 90846.25
$ws.Range("J137").Value = 5558.75
$ws.Range("K137").Value = 272538.75
$ws.Range("L137").Value = 16676.25
$ws.Range("M137").Value = -269988.75
$ws.Range("N137").Value = -21776.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5568.721
$ws.Range("I32").Value = 3686.4
$ws.Range("J32").Value = 9154.096
$ws.Range("K32").Value = 3686.4
$ws.Range("L32").Value = 9154.096
$ws.Range("M32").Value = -3399.4
$ws.Range("N32").Value = -9728.096
$ws.Range("H45").Value = 7995871
$ws.Range("I45").Value = 17983916
$ws.Range("K45").Value = 17983916
$ws.Range("M45").Value = -17983539
$ws.Range("H74").Value = 176217.27
$ws.Range("I74").Value = 81132.87
$ws.Range("J74").Value = 379969.56
$ws.Range("K74").Value = 81132.87
$ws.Range("L74").Value = 379969.56
$ws.Range("M74").Value = -80258.87
$ws.Range("N74").Value = -381717.56
$ws.Range("H77").Value = 176217.27
$ws.Range("I77").Value = 81132.87
$ws.Range("J77").Value = 379969.56
$ws.Range("K77").Value = 405664.35
$ws.Range("L77").Value = 1899847.8
$ws.Range("M77").Value = -401296.35
$ws.Range("N77").Value = -1908583.8
$ws.Range("H97").Value = 986458.25
$ws.Range("I97").Value = 1118865.5
$ws.Range("K97").Value = 1118865.5
$ws.Range("M97").Value = -1118369.5
$ws.Range("H110").Value = 1323504.9
$ws.Range("I110").Value = 1736969.6
$ws.Range("K110").Value = 1736969.6
$ws.Range("M110").Value = -1734924.6
$ws.Range("H122").Value = 835613.4
$ws.Range("I122").Value = 2296.389
$ws.Range("J122").Value = 2978428.5
$ws.Range("K122").Value = 6889.167
$ws.Range("L122").Value = 8935285.5
$ws.Range("M122").Value = -4439.167
$ws.Range("N122").Value = -8940185.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2623.7334
$ws.Range("I20").Value = 2349.5557
$ws.Range("J20").Value = 3035
$ws.Range("K20").Value = 2349.5557
$ws.Range("L20").Value = 3035
$ws.Range("M20").Value = -2102.5557
$ws.Range("N20").Value = -3529
$ws.Range("H99").Value = 7574433
$ws.Range("I99").Value = 14388113
$ws.Range("K99").Value = 14388113
$ws.Range("M99").Value = -14386615
$ws.Range("H107").Value = 8930088
$ws.Range("I107").Value = 8930088
$ws.Range("K107").Value = 8930088
$ws.Range("M107").Value = -8928168
$ws.Range("H134").Value = 4971.1763
$ws.Range("I134").Value = 1071.8
$ws.Range("J134").Value = 10541.714
$ws.Range("K134").Value = 3215.4
$ws.Range("L134").Value = 31625.142
$ws.Range("M134").Value = -680.3999999999996
$ws.Range("N134").Value = -36695.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1525.5
$ws.Range("I19").Value = 50
$ws.Range("K19").Value = 50
$ws.Range("M19").Value = 120
$ws.Range("H24").Value = 1525.5
$ws.Range("I24").Value = 50
$ws.Range("K24").Value = 50
$ws.Range("M24").Value = 120
$ws.Range("H58").Value = 4166.5835
$ws.Range("I58").Value = 3998.5715
$ws.Range("J58").Value = 4401.8
$ws.Range("K58").Value = 3998.5715
$ws.Range("L58").Value = 4401.8
$ws.Range("M58").Value = -3795.5715
$ws.Range("N58").Value = -4807.8
$ws.Range("H99").Value = 3520.8235
$ws.Range("I99").Value = 2599.3333
$ws.Range("J99").Value = 4557.5
$ws.Range("K99").Value = 2599.3333
$ws.Range("L99").Value = 4557.5
$ws.Range("M99").Value = -1101.3333
$ws.Range("N99").Value = -7553.5
$ws.Range("H109").Value = 48595.4
$ws.Range("J109").Value = 48595.4
$ws.Range("L109").Value = 48595.4
$ws.Range("N109").Value = -50675.4
$ws.Range("H126").Value = 3520.8235
$ws.Range("I126").Value = 2599.3333
$ws.Range("J126").Value = 4557.5
$ws.Range("K126").Value = 7797.999899999999
$ws.Range("L126").Value = 13672.5
$ws.Range("M126").Value = -5327.999899999999
$ws.Range("N126").Value = -18612.5
$ws.Range("H132").Value = 96368.5
$ws.Range("I132").Value = 68232.2
$ws.Range("K132").Value = 204696.6
$ws.Range("M132").Value = -202166.6
$ws.Range("H136").Value = 4166.5835
$ws.Range("I136").Value = 3998.5715
$ws.Range("J136").Value = 4401.8
$ws.Range("K136").Value = 11995.7145
$ws.Range("L136").Value = 13205.4
$ws.Range("M136").Value = -9445.7145
$ws.Range("N136").Value = -18305.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 93.60869599999999
$ws.Range("J38").Value = 94.5
$ws.Range("L38").Value = 283.5
$ws.Range("N38").Value = -977.5
$ws.Range("H56").Value = 10422778
$ws.Range("I56").Value = 10422778
$ws.Range("K56").Value = 10422778
$ws.Range("M56").Value = -10422248
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 2740.8572
$ws.Range("J132").Value = 2884.1428
$ws.Range("L132").Value = 25957.2852
$ws.Range("N132").Value = -31017.2852
$ws.Range("H139").Value = 71430940
$ws.Range("I139").Value = 125002030
$ws.Range("J139").Value = 2833
$ws.Range("K139").Value = 375006090
$ws.Range("L139").Value = 8499
$ws.Range("M139").Value = -375000950
$ws.Range("N139").Value = -18779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1879504.2
$ws.Range("I80").Value = 3050449.8
$ws.Range("J80").Value = 5991.4
$ws.Range("K80").Value = 3050449.8
$ws.Range("L80").Value = 5991.4
$ws.Range("M80").Value = -3049451.8
$ws.Range("N80").Value = -7987.4
$ws.Range("H83").Value = 1879504.2
$ws.Range("I83").Value = 3050449.8
$ws.Range("J83").Value = 5991.4
$ws.Range("K83").Value = 15252249
$ws.Range("L83").Value = 29957
$ws.Range("M83").Value = -15247257
$ws.Range("N83").Value = -39941
$ws.Range("H102").Value = 13595460
$ws.Range("I102").Value = 27780026
$ws.Range("K102").Value = 27780026
$ws.Range("M102").Value = -27778404
$ws.Range("H132").Value = 4077.0625
$ws.Range("I132").Value = 3463.5
$ws.Range("J132").Value = 5099.6665
$ws.Range("K132").Value = 10390.5
$ws.Range("L132").Value = 15298.9995
$ws.Range("M132").Value = -7860.5
$ws.Range("N132").Value = -20358.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 74767.414
$ws.Range("I22").Value = 148684.83
$ws.Range("J22").Value = 850
$ws.Range("K22").Value = 148684.83
$ws.Range("L22").Value = 850
$ws.Range("M22").Value = -148389.83
$ws.Range("N22").Value = -1440
$ws.Range("H27").Value = 74767.414
$ws.Range("I27").Value = 148684.83
$ws.Range("J27").Value = 850
$ws.Range("K27").Value = 148684.83
$ws.Range("L27").Value = 850
$ws.Range("M27").Value = -148577.83
$ws.Range("N27").Value = -1064
$ws.Range("H55").Value = 2087.3142
$ws.Range("I55").Value = 1627.8
$ws.Range("J55").Value = 3236.1
$ws.Range("K55").Value = 1627.8
$ws.Range("L55").Value = 3236.1
$ws.Range("M55").Value = -1454.8
$ws.Range("N55").Value = -3582.1
$ws.Range("H61").Value = 5559085
$ws.Range("I61").Value = 6539589
$ws.Range("K61").Value = 6539589
$ws.Range("M61").Value = -6539387
$ws.Range("H68").Value = 2543.0833
$ws.Range("I68").Value = 2494.6667
$ws.Range("J68").Value = 2688.3333
$ws.Range("K68").Value = 2494.6667
$ws.Range("L68").Value = 2688.3333
$ws.Range("M68").Value = -1745.6667
$ws.Range("N68").Value = -4186.3333
$ws.Range("H71").Value = 2543.0833
$ws.Range("I71").Value = 2494.6667
$ws.Range("J71").Value = 2688.3333
$ws.Range("K71").Value = 12473.3335
$ws.Range("L71").Value = 13441.6665
$ws.Range("M71").Value = -8729.333500000001
$ws.Range("N71").Value = -20929.6665
$ws.Range("H113").Value = 5559085
$ws.Range("I113").Value = 6539589
$ws.Range("K113").Value = 6539589
$ws.Range("M113").Value = -6537419
$ws.Range("H132").Value = 5771.8213
$ws.Range("I132").Value = 6080.7295
$ws.Range("K132").Value = 18242.1885
$ws.Range("M132").Value = -15712.1885
$ws.Range("H136").Value = 47857.49
$ws.Range("I136").Value = 72260.2
$ws.Range("J136").Value = 3627.5625
$ws.Range("K136").Value = 216780.6
$ws.Range("L136").Value = 10882.6875
$ws.Range("M136").Value = -214230.6
$ws.Range("N136").Value = -15982.6875
